# Replace en-dash (U+2013) with a plain hyphen-minus (U+002D) in the
# opening-hours strings found in column D, rows 5 through 19, of the
# "Toiletten" worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Toiletten")

$enDash = [char]0x2013

for ($row = 5; $row -le 19; $row++) {
    $cell = $ws.Cells.Item($row, 4)  # column D
    $value = $cell.Value2
    if ($value -ne $null -and $value.ToString().Contains($enDash)) {
        $cell.Value2 = $value.ToString().Replace($enDash, "-")
    }
}

# Update the active selection on the sheet from E18 to E19.
$ws.Range("E19").Select()
